$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A65").Value = 'https://web3.career/lead-web3-qa-engineer-dia/53252'
$ws.Range("A66").Value = 'https://web3.career/head-of-engineering-cere-network/54577'
$ws.Range("A67").Value = 'https://web3.career/senior-web3-front-end-developer-invarch-network/53229'
$ws.Range("A68").Value = 'https://web3.career/growth-manager-everest-ventures-group-evg/54555'
$ws.Range("A69").Value = 'https://web3.career/data-engineer-trilitech/54603'
$ws.Range("A70").Value = 'https://web3.career/events-manager-trilitech/54602'
$ws.Range("A71").Value = 'https://web3.career/legal-counsel-regulatory-compliance-swissborg/54601'
$ws.Range("A72").Value = 'https://web3.career/director-of-sales-hong-kong-bitgo/54600'
$ws.Range("A73").Value = 'https://web3.career/senior-developer-relations-engineer-web3-monadlabs/54599'
$ws.Range("A74").Value = 'https://web3.career/fp-a-intern-fall-2023-shiftmarkets/54598'
$ws.Range("A75").Value = 'https://web3.career/defi-community-manager-shiftmarkets/54597'
$ws.Range("A76").Value = 'https://web3.career/solutions-architect-axelarnetwork/54596'
$ws.Range("A77").Value = 'https://web3.career/freelance-translator-traditional-chinese-product-localization-crypto-com/54595'
$ws.Range("A78").Value = 'https://web3.career/freelance-translator-japanese-product-localization-crypto-com/54594'
$ws.Range("A79").Value = 'https://web3.career/freelance-translator-indonesian-product-localization-crypto-com/54593'
$ws.Range("A80").Value = 'https://web3.career/freelance-translator-danish-product-localization-crypto-com/54592'
$ws.Range("A81").Value = 'https://web3.career/content-marketer-demand-generation-uniswaplabs/54591'
$ws.Range("A82").Value = 'https://web3.career/marketing-lead-with-foundation/54590'
$ws.Range("A83").Value = 'https://web3.career/finance-accounting-operations-staff-pintu/54589'
$ws.Range("A84").Value = 'https://web3.career/site-reliability-engineer-pintu/54588'
$ws.Range("A85").Value = 'https://web3.career/senior-network-engineer-ripple/54587'
$ws.Range("A86").Value = 'https://web3.career/international-tax-director-ripple/54586'
$ws.Range("A87").Value = 'https://web3.career/partner-6-executive-assistant-crypto-a16z/54585'
$ws.Range("A88").Value = 'https://web3.career/analytics-engineer-gemini/54584'
$ws.Range("A89").Value = 'https://web3.career/group-product-manager-base-ecosystem-base/54583'
$ws.Range("A90").Value = 'https://web3.career/commercial-regulatory-legal-counsel-gnosis/54582'
$ws.Range("A91").Value = 'https://web3.career/senior-lead-golang-engineer-cere-network/52169'
$ws.Range("A92").Value = 'https://web3.career/lead-software-engineer-blockchain-cere-network/52168'
$ws.Range("A93").Value = 'https://web3.career/lead-software-engineer-blockchain-cere-network/54749'
$ws.Range("A94").Value = 'https://web3.career/head-of-cybersecurity-coinclan-ou/54694'
$ws.Range("A95").Value = 'https://web3.career/senior-software-engineer-java-prime-services-okx/54761'
$ws.Range("A96").Value = 'https://web3.career/senior-finance-manager-crypto-com/54760'
$ws.Range("A97").Value = 'https://web3.career/senior-product-manager-pay-crypto-com/54759'
$ws.Range("A98").Value = 'https://web3.career/product-director-operations-platform-okx/54758'
$ws.Range("A99").Value = 'https://web3.career/marketing-designer-huddle01/54757'
$ws.Range("A100").Value = 'https://web3.career/senior-product-strategy-manager-growth-okx/54756'
$ws.Range("A101").Value = 'https://web3.career/senior-product-manager-risk-okx/54755'
$ws.Range("A102").Value = 'https://web3.career/product-operations-lead-growth-platform-okx/54754'
$ws.Range("A103").Value = 'https://web3.career/senior-product-manager-asset-trading-immutable/54753'
$ws.Range("A104").Value = 'https://web3.career/binance-accelerator-program-angel-program-coordinator-binance/54752'
$ws.Range("A105").Value = 'https://web3.career/fp-a-manager-chainlink/54751'
$ws.Range("A106").Value = 'https://web3.career/direktor-fur-projekte-und-produkte-zaubar/54750'
$ws.Range("A107").Value = 'https://web3.career/senior-growth-lead-keyp/53471'
$ws.Range("A108").Value = 'https://web3.career/managing-director-partner-theflowerscompany/54748'
$ws.Range("A109").Value = 'https://web3.career/vice-president-dittopr/54746'
$ws.Range("A110").Value = 'https://web3.career/software-development-engineer-in-test-exchange-pintu/54745'
$ws.Range("A111").Value = 'https://web3.career/senior-manager-people-technology-ripple/54744'
$ws.Range("A112").Value = 'https://web3.career/senior-full-stack-software-engineer-ripple/54743'
$ws.Range("A113").Value = 'https://web3.career/staff-cloud-infrastructure-engineer-gemini/54741'
$ws.Range("A114").Value = 'https://web3.career/data-engineer-consumer-experience-gemini/54740'
$ws.Range("A115").Value = 'https://web3.career/director-of-marketing-nansen/54762'
